$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spending data")
$ws.Columns.Item(2).Insert()
